$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell value updates (terrain grid edits) ---
$ws.Range("BK55").Value = 1
$ws.Range("BL55").Value = 1
$ws.Range("BM55").Value = 1
$ws.Range("BN55").Value = 1
$ws.Range("BO55").Value = 1
$ws.Range("BK56").Value = 1
$ws.Range("BL56").Value = 1
$ws.Range("BM56").Value = 1
$ws.Range("BN56").Value = 1
$ws.Range("BO56").Value = 1
$ws.Range("BK57").Value = 1
$ws.Range("BL57").Value = 1
$ws.Range("BM57").Value = 1
$ws.Range("BN57").Value = 1
$ws.Range("BO57").Value = 1
$ws.Range("BH58").Value = 1
$ws.Range("BI58").Value = 1
$ws.Range("BJ58").Value = 1
$ws.Range("BK58").Value = 1
$ws.Range("BL58").Value = 1
$ws.Range("BM58").Value = 1
$ws.Range("BN58").Value = 1
$ws.Range("BH59").Value = 1
$ws.Range("BI59").Value = 1
$ws.Range("BJ59").Value = 1
$ws.Range("BK59").Value = 1
$ws.Range("BL59").Value = 1
$ws.Range("BM59").Value = 1
$ws.Range("BN59").Value = 1
$ws.Range("BR59").Value = 1
$ws.Range("BS59").Value = 1
$ws.Range("BT59").Value = 1
$ws.Range("BU59").Value = 1
$ws.Range("BV59").Value = 1
$ws.Range("BJ60").Value = 1
$ws.Range("BK60").Value = 1
$ws.Range("BL60").Value = 1
$ws.Range("BM60").Value = 1
$ws.Range("BN60").Value = 1
$ws.Range("BR60").Value = 1
$ws.Range("BS60").Value = 1
$ws.Range("BT60").Value = 1
$ws.Range("BU60").Value = 1
$ws.Range("BV60").Value = 1
$ws.Range("BR61").Value = 1
$ws.Range("BS61").Value = 1
$ws.Range("BT61").Value = 1
$ws.Range("BU61").Value = 1
$ws.Range("BV61").Value = 1
$ws.Range("BR62").Value = 1
$ws.Range("BS62").Value = 1
$ws.Range("BT62").Value = 1
$ws.Range("BU62").Value = 1
$ws.Range("BV62").Value = 1
$ws.Range("P62").Value = 3
$ws.Range("Q62").Value = 3
$ws.Range("S62").Value = 3
$ws.Range("T62").Value = 3
$ws.Range("U62").Value = 3
$ws.Range("CE68").Value = 3
$ws.Range("CF68").Value = 3
$ws.Range("CG68").Value = 3
$ws.Range("CH68").Value = 3
$ws.Range("CI68").Value = 3
$ws.Range("CJ68").Value = 4
$ws.Range("CK68").Value = 4
$ws.Range("CL68").Value = 3
$ws.Range("CM68").Value = 3
$ws.Range("CN68").Value = 3
$ws.Range("CJ69").Value = 4
$ws.Range("CK69").Value = 4
$ws.Range("CJ70").Value = 4
$ws.Range("CK70").Value = 4
$ws.Range("CA71").Value = 4
$ws.Range("CB71").Value = 4
$ws.Range("CC71").Value = 4
$ws.Range("CD71").Value = 4
$ws.Range("CE71").Value = 4
$ws.Range("CF71").Value = 4
$ws.Range("CG71").Value = 4
$ws.Range("CH71").Value = 4
$ws.Range("CI71").Value = 4
$ws.Range("CJ71").Value = 4
$ws.Range("CK71").Value = 4
$ws.Range("CA72").Value = 4
$ws.Range("CB72").Value = 4
$ws.Range("CC72").Value = 4
$ws.Range("CD72").Value = 4
$ws.Range("CE72").Value = 4
$ws.Range("CF72").Value = 4
$ws.Range("CG72").Value = 4
$ws.Range("CH72").Value = 4
$ws.Range("CI72").Value = 4
$ws.Range("CJ72").Value = 4
$ws.Range("CK72").Value = 4
$ws.Range("CA73").Value = 4
$ws.Range("CB73").Value = 4
$ws.Range("CA74").Value = 4
$ws.Range("CB74").Value = 4
$ws.Range("BO75").Value = 4
$ws.Range("BP75").Value = 4
$ws.Range("BQ75").Value = 4
$ws.Range("BR75").Value = 4
$ws.Range("BS75").Value = 4
$ws.Range("BT75").Value = 4
$ws.Range("BU75").Value = 4
$ws.Range("BV75").Value = 4
$ws.Range("BW75").Value = 4
$ws.Range("BX75").Value = 4
$ws.Range("BY75").Value = 4
$ws.Range("BZ75").Value = 4
$ws.Range("CA75").Value = 4
$ws.Range("CB75").Value = 4
$ws.Range("BO76").Value = 4
$ws.Range("BP76").Value = 4
$ws.Range("BQ76").Value = 4
$ws.Range("BR76").Value = 4
$ws.Range("BS76").Value = 4
$ws.Range("BT76").Value = 4
$ws.Range("BU76").Value = 4
$ws.Range("BV76").Value = 4
$ws.Range("BW76").Value = 4
$ws.Range("BX76").Value = 4
$ws.Range("BY76").Value = 4
$ws.Range("BZ76").Value = 4
$ws.Range("CA76").Value = 4
$ws.Range("CB76").Value = 4
$ws.Range("BO77").Value = 3
$ws.Range("BP77").Value = 3
$ws.Range("CA77").Value = 4
$ws.Range("CB77").Value = 4
$ws.Range("CA78").Value = 4
$ws.Range("CB78").Value = 4
$ws.Range("CA79").Value = 4
$ws.Range("CB79").Value = 4
$ws.Range("BO80").Value = 3
$ws.Range("CA80").Value = 4
$ws.Range("CB80").Value = 4
$ws.Range("BO81").Value = 3
$ws.Range("CA81").Value = 4
$ws.Range("CB81").Value = 4
$ws.Range("CA82").Value = 4
$ws.Range("CB82").Value = 4
$ws.Range("O82").Value = 4
$ws.Range("P82").Value = 4
$ws.Range("Q82").Value = 4
$ws.Range("R82").Value = 4
$ws.Range("S82").Value = 4
$ws.Range("T82").Value = 4
$ws.Range("U82").Value = 4
$ws.Range("V82").Value = 4
$ws.Range("W82").Value = 4
$ws.Range("X82").Value = 4
$ws.Range("CA83").Value = 4
$ws.Range("CB83").Value = 4
$ws.Range("O83").Value = 4
$ws.Range("P83").Value = 4
$ws.Range("Q83").Value = 4
$ws.Range("R83").Value = 4
$ws.Range("S83").Value = 4
$ws.Range("T83").Value = 4
$ws.Range("U83").Value = 4
$ws.Range("V83").Value = 4
$ws.Range("W83").Value = 4
$ws.Range("X83").Value = 4
$ws.Range("CA84").Value = 4
$ws.Range("CB84").Value = 4
$ws.Range("O84").Value = 4
$ws.Range("P84").Value = 4
$ws.Range("CA85").Value = 4
$ws.Range("CB85").Value = 4
$ws.Range("O85").Value = 4
$ws.Range("P85").Value = 4
$ws.Range("R85").Value = 3
$ws.Range("T85").Value = 3
$ws.Range("U85").Value = 3
$ws.Range("V85").Value = 3
$ws.Range("W85").Value = 3
$ws.Range("CA86").Value = 4
$ws.Range("CB86").Value = 4
$ws.Range("O86").Value = 4
$ws.Range("P86").Value = 4
$ws.Range("R86").Value = 3
$ws.Range("W86").Value = 3
$ws.Range("CA87").Value = 4
$ws.Range("CB87").Value = 4
$ws.Range("O87").Value = 4
$ws.Range("P87").Value = 4
$ws.Range("R87").Value = 3
$ws.Range("W87").Value = 3
$ws.Range("CA88").Value = 4
$ws.Range("CB88").Value = 4
$ws.Range("O88").Value = 4
$ws.Range("P88").Value = 4
$ws.Range("R88").Value = 3
$ws.Range("S88").Value = 3
$ws.Range("T88").Value = 3
$ws.Range("U88").Value = 3
$ws.Range("V88").Value = 3
$ws.Range("W88").Value = 3
$ws.Range("CA89").Value = 4
$ws.Range("CB89").Value = 4
$ws.Range("O89").Value = 4
$ws.Range("P89").Value = 4
$ws.Range("BO90").Value = 3
$ws.Range("CA90").Value = 4
$ws.Range("CB90").Value = 4
$ws.Range("CC90").Value = 3
$ws.Range("CD90").Value = 3
$ws.Range("O90").Value = 4
$ws.Range("P90").Value = 4
$ws.Range("Q90").Value = 4
$ws.Range("R90").Value = 4
$ws.Range("S90").Value = 4
$ws.Range("T90").Value = 4
$ws.Range("U90").Value = 4
$ws.Range("V90").Value = 4
$ws.Range("W90").Value = 4
$ws.Range("X90").Value = 4
$ws.Range("BO91").Value = 3
$ws.Range("BQ91").Value = 3
$ws.Range("BR91").Value = 3
$ws.Range("BS91").Value = 3
$ws.Range("BT91").Value = 3
$ws.Range("BU91").Value = 3
$ws.Range("CA91").Value = 4
$ws.Range("CB91").Value = 4
$ws.Range("CC91").Value = 3
$ws.Range("CD91").Value = 3
$ws.Range("O91").Value = 4
$ws.Range("P91").Value = 4
$ws.Range("Q91").Value = 4
$ws.Range("R91").Value = 4
$ws.Range("S91").Value = 4
$ws.Range("T91").Value = 4
$ws.Range("U91").Value = 4
$ws.Range("V91").Value = 4
$ws.Range("W91").Value = 4
$ws.Range("X91").Value = 4

# --- View state: selection + scroll position ---
$ws.Range("CL68:CN68").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 25
$win.ScrollColumn = 1
